# Auto-generated edit script
# Updates betting-odds values in Sheet1 per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$row2 = @{
    "G2" = 3.3
    "I2" = 2.55
    "J2" = 4.33
    "K2" = 1.8
    "L2" = 3.5
    "M2" = 1.14
    "N2" = 5.5
    "O2" = 1.67
    "P2" = 2.1
    "Q2" = 3.4
    "R2" = 1.33
    "S2" = 1.75
    "T2" = 2.05
    "W2" = 6.5
    "X2" = 13
    "Z2" = 41
    "AA2" = 41
    "AC2" = 5
    "AH2" = 10
    "AJ2" = 26
    "AK2" = 29
    "AN2" = 5
    "AO2" = 23
    "AX2" = 17
}
foreach ($addr in $row2.Keys) {
    $ws.Range($addr).Value = $row2[$addr]
}

# Row 7
$row7 = @{
    "G7" = 2.1
    "I7" = 3.8
    "J7" = 2.88
    "K7" = 2.05
    "L7" = 4.33
    "U7" = 1.91
    "V7" = 1.8
    "W7" = 6.5
    "X7" = 9.5
    "Y7" = 9.5
    "Z7" = 19
    "AE7" = 15
    "AG7" = 9.5
    "AI7" = 13
    "AK7" = 34
    "AM7" = 351
    "AO7" = 12
    "AP7" = 23
    "AU7" = 8.5
    "AX7" = 21
    "AZ7" = 67
    "BA7" = 101
    "BB7" = 251
}
foreach ($addr in $row7.Keys) {
    $ws.Range($addr).Value = $row7[$addr]
}

# Row 8
$row8 = @{
    "G8" = 1.53
    "H8" = 3.7
    "I8" = 7
    "J8" = 2.2
    "K8" = 2.1
    "O8" = 1.36
    "P8" = 3
    "Q8" = 2.15
    "R8" = 1.67
    "S8" = 1.44
    "T8" = 2.63
    "Z8" = 10
    "AC8" = 8
    "AD8" = 7.5
    "AG8" = 15
    "AH8" = 34
    "AJ8" = 81
    "AO8" = 8
    "AQ8" = 26
    "AT8" = 2.63
    "AU8" = 10
    "AW8" = 7.5
}
foreach ($addr in $row8.Keys) {
    $ws.Range($addr).Value = $row8[$addr]
}

# Row 11
$row11 = @{
    "G11" = 2.75
    "H11" = 3.4
    "I11" = 2.45
    "J11" = 3.6
    "L11" = 3.25
    "M11" = 1.07
    "N11" = 9
    "O11" = 1.4
    "P11" = 2.75
    "Q11" = 2.25
    "R11" = 1.62
    "W11" = 7.5
    "X11" = 13
    "Y11" = 11
    "Z11" = 29
    "AA11" = 26
    "AC11" = 8
    "AE11" = 17
    "AG11" = 7
    "AH11" = 11
    "AI11" = 10
    "AJ11" = 23
    "AK11" = 21
    "AN11" = 4.75
    "AO11" = 17
    "AP11" = 29
    "AQ11" = 51
    "AR11" = 81
    "AW11" = 4.5
    "AX11" = 15
    "AY11" = 26
    "AZ11" = 51
    "BA11" = 81
}
foreach ($addr in $row11.Keys) {
    $ws.Range($addr).Value = $row11[$addr]
}

# Row 12
$row12 = @{
    "M12" = 1.07
    "N12" = 8.5
    "X12" = 23
    "AN12" = 6.5
    "AQ12" = 101
    "AX12" = 9.5
}
foreach ($addr in $row12.Keys) {
    $ws.Range($addr).Value = $row12[$addr]
}

# Row 14
$row14 = @{
    "G14" = 2.15
    "H14" = 2.9
    "I14" = 3.9
    "J14" = 3.1
    "L14" = 5
    "M14" = 1.13
    "N14" = 6
    "O14" = 1.62
    "P14" = 2.2
    "S14" = 1.67
    "T14" = 2.1
    "W14" = 5
    "X14" = 8.5
    "Y14" = 11
    "Z14" = 19
    "AA14" = 23
    "AB14" = 41
    "AE14" = 23
    "AG14" = 7.5
    "AH14" = 17
    "AI14" = 15
    "AJ14" = 41
    "AK14" = 41
    "AN14" = 4
    "AO14" = 13
    "AQ14" = 51
    "AR14" = 101
    "AT14" = 2.1
    "AW14" = 5.5
    "AX14" = 26
    "AZ14" = 101
    "BA14" = 151
}
foreach ($addr in $row14.Keys) {
    $ws.Range($addr).Value = $row14[$addr]
}

# Row 15
$row15 = @{
    "G15" = 1.4
    "I15" = 8
    "K15" = 2.2
    "L15" = 9
    "AH15" = 41
    "AI15" = 26
    "AL15" = 81
    "AP15" = 26
    "AQ15" = 21
    "AW15" = 9.5
    "AX15" = 51
    "AZ15" = 301
    "BA15" = 351
}
foreach ($addr in $row15.Keys) {
    $ws.Range($addr).Value = $row15[$addr]
}

Write-Host "Applied all odds updates."